$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 476
$ws.Range("J58").Value = 983.5
$ws.Range("L58").Value = 2950.5
$ws.Range("N58").Value = -3250.5
$ws.Range("H70").Value = 1461.5385
$ws.Range("I70").Value = 1462
$ws.Range("J70").Value = 1460
$ws.Range("K70").Value = 4386
$ws.Range("L70").Value = 4380
$ws.Range("M70").Value = -4116
$ws.Range("N70").Value = -4920
$ws.Range("H73").Value = 1461.5385
$ws.Range("I73").Value = 1462
$ws.Range("J73").Value = 1460
$ws.Range("K73").Value = 4386
$ws.Range("L73").Value = 4380
$ws.Range("M73").Value = -3450
$ws.Range("N73").Value = -6252
$ws.Range("H86").Value = 17654.166
$ws.Range("I86").Value = 933
$ws.Range("J86").Value = 34375.332
$ws.Range("K86").Value = 933
$ws.Range("L86").Value = 34375.332
$ws.Range("M86").Value = 190
$ws.Range("N86").Value = -36621.332
$ws.Range("H89").Value = 17654.166
$ws.Range("I89").Value = 933
$ws.Range("J89").Value = 34375.332
$ws.Range("K89").Value = 4665
$ws.Range("L89").Value = 171876.66
$ws.Range("M89").Value = 951
$ws.Range("N89").Value = -183108.66
$ws.Range("H129").Value = 334392.5
$ws.Range("J129").Value = 401199
$ws.Range("L129").Value = 1203597
$ws.Range("N129").Value = -1213597
$ws.Range("H131").Value = 2235.9565
$ws.Range("I131").Value = 1589.1875
$ws.Range("J131").Value = 3714.2856
$ws.Range("K131").Value = 4767.5625
$ws.Range("L131").Value = 11142.8568
$ws.Range("M131").Value = 272.4375
$ws.Range("N131").Value = -21222.8568
$ws.Range("H138").Value = 4660.6084
$ws.Range("J138").Value = 4659.989
$ws.Range("L138").Value = 13979.967
$ws.Range("N138").Value = -24259.967

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2758
$ws.Range("I2").Value = 2758
$ws.Range("K2").Value = 2758
$ws.Range("M2").Value = -2645
$ws.Range("H32").Value = 12545.167
$ws.Range("I32").Value = 10363.3125
$ws.Range("K32").Value = 10363.3125
$ws.Range("M32").Value = -10076.3125
$ws.Range("H63").Value = 15626150
$ws.Range("I63").Value = 2300
$ws.Range("J63").Value = 31250000
$ws.Range("K63").Value = 2300
$ws.Range("L63").Value = 31250000
$ws.Range("M63").Value = -1614
$ws.Range("N63").Value = -31251372
$ws.Range("H66").Value = 15626150
$ws.Range("I66").Value = 2300
$ws.Range("J66").Value = 31250000
$ws.Range("K66").Value = 11500
$ws.Range("L66").Value = 156250000
$ws.Range("M66").Value = -8068
$ws.Range("N66").Value = -156256864
$ws.Range("H116").Value = 2758
$ws.Range("I116").Value = 2758
$ws.Range("K116").Value = 2758
$ws.Range("M116").Value = -464

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2758
$ws.Range("I3").Value = 2758
$ws.Range("K3").Value = 2758
$ws.Range("M3").Value = -2644
$ws.Range("H20").Value = 4900
$ws.Range("I20").Value = 4716.6665
$ws.Range("J20").Value = 6000
$ws.Range("K20").Value = 4716.6665
$ws.Range("L20").Value = 6000
$ws.Range("M20").Value = -4469.6665
$ws.Range("N20").Value = -6494

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1508.8462
$ws.Range("I5").Value = 1210.12
$ws.Range("J5").Value = 2042.2858
$ws.Range("K5").Value = 3630.36
$ws.Range("L5").Value = 6126.857400000001
$ws.Range("M5").Value = -3518.36
$ws.Range("N5").Value = -6350.857400000001
$ws.Range("H8").Value = 334
$ws.Range("I8").Value = 334
$ws.Range("K8").Value = 1002
$ws.Range("M8").Value = -863
$ws.Range("H33").Value = 62.857143
$ws.Range("H44").Value = 382.3889
$ws.Range("I44").Value = 341.85715
$ws.Range("J44").Value = 524.25
$ws.Range("K44").Value = 1025.57145
$ws.Range("L44").Value = 1572.75
$ws.Range("M44").Value = -627.5714499999999
$ws.Range("N44").Value = -2368.75
$ws.Range("H131").Value = 103873.875
$ws.Range("J131").Value = 113163.35
$ws.Range("L131").Value = 339490.05
$ws.Range("N131").Value = -349570.05
$ws.Range("H135").Value = 1508.8462
$ws.Range("I135").Value = 1210.12
$ws.Range("J135").Value = 2042.2858
$ws.Range("K135").Value = 10891.08
$ws.Range("L135").Value = 18380.5722
$ws.Range("M135").Value = -8356.079999999998
$ws.Range("N135").Value = -23450.5722
$ws.Range("H137").Value = 30310976
$ws.Range("I137").Value = 3066.3333
$ws.Range("J137").Value = 41676444
$ws.Range("K137").Value = 9198.999899999999
$ws.Range("L137").Value = 125029332
$ws.Range("M137").Value = -4098.999899999999
$ws.Range("N137").Value = -125039532

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3947.1177
$ws.Range("I80").Value = 3633
$ws.Range("J80").Value = 4118.4546
$ws.Range("K80").Value = 3633
$ws.Range("L80").Value = 4118.4546
$ws.Range("M80").Value = -2635
$ws.Range("N80").Value = -6114.4546
$ws.Range("H83").Value = 3947.1177
$ws.Range("I83").Value = 3633
$ws.Range("J83").Value = 4118.4546
$ws.Range("K83").Value = 18165
$ws.Range("L83").Value = 20592.273
$ws.Range("M83").Value = -13173
$ws.Range("N83").Value = -30576.273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1334.5
$ws.Range("I46").Value = 1135.8
$ws.Range("J46").Value = 1392.9412
$ws.Range("K46").Value = 1135.8
$ws.Range("L46").Value = 1392.9412
$ws.Range("M46").Value = -947.8
$ws.Range("N46").Value = -1768.9412
$ws.Range("H68").Value = 2927.8572
$ws.Range("I68").Value = 2750.5
$ws.Range("J68").Value = 3164.3333
$ws.Range("K68").Value = 2750.5
$ws.Range("L68").Value = 3164.3333
$ws.Range("M68").Value = -2001.5
$ws.Range("N68").Value = -4662.3333
$ws.Range("H71").Value = 2927.8572
$ws.Range("I71").Value = 2750.5
$ws.Range("J71").Value = 3164.3333
$ws.Range("K71").Value = 13752.5
$ws.Range("L71").Value = 15821.6665
$ws.Range("M71").Value = -10008.5
$ws.Range("N71").Value = -23309.6665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2316.5557
$ws.Range("H84").Value = 2316.5557
$ws.Range("H107").Value = 101.5
$ws.Range("I107").Value = 101.5
$ws.Range("K107").Value = 304.5
$ws.Range("M107").Value = 1615.5
